# Apply "Add data for 2022-05-13" update:
# - Rename sheet "Through 2022-05-04" -> "Through 2022-05-05"
# - Update header label "2022 (through 05-04)" -> "2022 (through 05-05)"
# - Update May 2022 value (I6) 12 -> 17
# - Update Total 2022 value (I14) 563 -> 568

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-05-05"

$ws.Range("I1").Value = "2022 (through 05-05)"

$ws.Range("I6").Value = 17

$ws.Range("I14").Value = 568
